# "predictions updated using weather data"
# Append 9 new rows (71-79) to Sheet1 continuing the "2021-01-09" prediction
# block (model group A=14) with updated Prediction values.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Column A on these rows must stay TEXT (e.g. "2021-01-09"), matching the
# existing cells above them - otherwise Excel auto-detects the string as a
# date and stores a serial number instead. Force the whole new block to
# Text format first, enter the values, then drop back to the default
# ("Normal") cell style so no stray formatting is left behind.
$aRange = $ws.Range("A71:A79")
$aRange.NumberFormat = "@"

$data = @(
    @{ Row = 71; A = "2021-01-09"; B = "10 Jan -- 16 Jan 2021"; C = 3333.57; D = 935.4;   E = 2398.17; F = "KNN"; G = 0.8;  H = 61.5;  I = 1952.14; J = 2186.81; K = 70.96 },
    @{ Row = 72; A = "2021-01-09"; B = "17 Jan -- 23 Jan 2021"; D = 1074.3;  F = "KNN" },
    @{ Row = 73; A = "2021-01-09"; B = "24 Jan -- 30 Jan 2021"; D = 811.55;  F = "KNN" },
    @{ Row = 74; A = "2021-01-09"; B = "31 Jan -- 06 Feb 2021"; D = 1419.86; F = "KNN" },
    @{ Row = 75; A = "2021-01-09"; B = "07 Feb -- 13 Feb 2021"; D = 1660.02; F = "KNN" },
    @{ Row = 76; A = "2021-01-09"; B = "14 Feb -- 20 Feb 2021"; D = 1516.25; F = "KNN" },
    @{ Row = 77; A = "2021-01-09"; B = "21 Feb -- 27 Feb 2021"; D = 528.57;  F = "KNN" },
    @{ Row = 78; A = "2021-01-09"; B = "28 Feb -- 06 Mar 2021"; D = 609.45;  F = "KNN" },
    @{ Row = 79; A = "2021-01-09"; B = "07 Mar -- 13 Mar 2021"; D = 1471.75; F = "KNN" }
)

foreach ($row in $data) {
    $r = $row.Row
    $ws.Cells.Item($r, 1).Value = $row.A
    $ws.Cells.Item($r, 2).Value = $row.B
    if ($row.ContainsKey("C")) { $ws.Cells.Item($r, 3).Value = $row.C }
    $ws.Cells.Item($r, 4).Value = $row.D
    if ($row.ContainsKey("E")) { $ws.Cells.Item($r, 5).Value = $row.E }
    $ws.Cells.Item($r, 6).Value = $row.F
    if ($row.ContainsKey("G")) { $ws.Cells.Item($r, 7).Value = $row.G }
    if ($row.ContainsKey("H")) { $ws.Cells.Item($r, 8).Value = $row.H }
    if ($row.ContainsKey("I")) { $ws.Cells.Item($r, 9).Value = $row.I }
    if ($row.ContainsKey("J")) { $ws.Cells.Item($r, 10).Value = $row.J }
    if ($row.ContainsKey("K")) { $ws.Cells.Item($r, 11).Value = $row.K }
}

# Restore the default style on column A for the new rows so only the
# number format (now back to General/shared-string text) differs from a
# freshly written cell - no leftover custom formatting.
$aRange.Style = "Normal"
